# Applies the commit "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
# - Adds 3 new worker rows to the account-statement table
# - Updates the totals (Valor Mora, Cant. Trabajadores, Cant. Periodos)
# - Reorders / edits a couple of the pre-existing data rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Update the summary header values
# ---------------------------------------------------------------------
$ws.Range("E11").Value = 176934      # VALOR MORA
$ws.Range("C13").Value = 5           # Cant. Trabajadores
$ws.Range("F13").Value = 4           # Cant. Periodos

# ---------------------------------------------------------------------
# 2) Make room in the worker table: it grows from 3 data rows (16-18)
#    to 6 data rows (16-21). Insert 3 fresh rows right above the
#    existing data block, then stamp them with the same formatting
#    used by the rest of the table (copied from the template row).
# ---------------------------------------------------------------------
$ws.Rows("16:18").Insert()

$ws.Range("B19:J19").Copy($ws.Range("B16:J16"))
$ws.Range("B19:J19").Copy($ws.Range("B17:J17"))
$ws.Range("B19:J19").Copy($ws.Range("B18:J18"))

# ---------------------------------------------------------------------
# 3) Fill in the final content of every worker row (16-21)
# ---------------------------------------------------------------------

# Row 16 - new worker
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1001969554"
$ws.Range("D16").Value = "EDER DE JESUS AVILA BERRIO"
$ws.Range("E16").Value = "2205"
$ws.Range("F16").Value = 40000
$ws.Range("G16").Value = 1000000

# Row 17 - previously the first data row; period updated
# (D17 is a purely-numeric label stored in a General-formatted column, so a
#  leading apostrophe is needed to keep it text, then the cell format is
#  restamped from a same-style neighbour to drop the quote-prefix marker)
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "ARISTIDES QUINTERO CHANTACA"
$ws.Range("D17").Value = "'2111"
$ws.Range("E17").Value = "2205"
$ws.Range("F17").Value = 16959
$ws.Range("G17").Value = 908526
$ws.Range("D16").Copy()
$ws.Range("D17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 18 - new worker
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1019120514"
$ws.Range("D18").Value = "JUAN FELIPE GUERRERO LONDOÑO"
$ws.Range("E18").Value = "2205"
$ws.Range("F18").Value = 40000
$ws.Range("G18").Value = 1000000

# Row 19 - new worker
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "19201423"
$ws.Range("D19").Value = "JORGE RICARDO GUERRERO MOLANO"
$ws.Range("E19").Value = "2205"
$ws.Range("F19").Value = 40000
$ws.Range("G19").Value = 1000000

# Row 20 - pre-existing worker/period, kept on the regular row style
$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "1005489479"
$ws.Range("D20").Value = "MARIA ELENA PEÑA CHAMORRO"
$ws.Range("E20").Value = "2201"
$ws.Range("F20").Value = 36341
$ws.Range("G20").Value = 1423500

# Row 21 - last table row (keeps the heavier bottom border of the
# original template's last row), pre-existing worker/period
$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "1005489479"
$ws.Range("D21").Value = "MARIA ELENA PEÑA CHAMORRO"
$ws.Range("E21").Value = "2112"
$ws.Range("F21").Value = 3634
$ws.Range("G21").Value = 1423500
